$wb = $excel.ActiveWorkbook

# Update the "想去人数" (interested count) column F on both the "展览"
# and "全部类型" sheets, which contain duplicated data.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 203
    $ws.Range("F3").Value = 154
    $ws.Range("F4").Value = 135
}
